# Changes to allow for unit testing:
# - Row 2 (SQL Server 2010) no longer references a MstrSkillID (clear D2,
#   and drop the MstrSkillID component from the generated INSERT formula).
# - Row 3 (Redditor) now points its MstrSkillID at row 2's SkillID.
# - Row 4 (4Chan Troll) now points its MstrSkillID at row 3's SkillID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove the MstrSkillID value entirely.
$ws.Range("D2").ClearContents()

# Row 2: rebuild the query formula without the MstrSkillID ($D$1 / RC[-3]) piece.
$ws.Range("G2").Formula = '=((((((((((((((((("INSERT INTO " & A2) &" (") & $B$1) & ",") & $C$1) & ",") & $E$1) & ",") & $F$1) & ") VALUES(''") & RC[-5]) &"'',''") & RC[-4]) &"'',''") & RC[-2]) &"'',''") & RC[-1]) & "'')"'

# Row 3: MstrSkillID now points at the SkillID used in row 2.
$ws.Range("D3").Value = "153B9006-1488-4CA1-950A-32E28F7B699D"

# Row 4: MstrSkillID now points at the SkillID used in row 3.
$ws.Range("D4").Value = "C87F23E9-8F8C-406D-9FBF-E15043179F09"
